$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the worksheet
$ws.Name = "Madera_5.0_4.0_54.0"

# Header info: material and dimensions
$ws.Range("B1").Value = "Madera"
$ws.Range("B2").Value = 5
$ws.Range("B3").Value = 4
$ws.Range("B4").Value = 54

# Data rows B7:AF10 (CREMER, SHARP, DAVY, ISO 12354-1:2001 results)
$row7 = @(9.93, 11.86, 13.87, 15.95, 17.89, 19.89, 21.97, 23.91, 25.84, 27.99, 29.93, 31.86, 7.86, 14.84, 19.23, 23.11, 26.78, 30.03, 33.17, 36.56, 39.58, 42.56, 45.62, 48.76, 51.69, 54.7, 57.82, 60.72, 63.62, 66.83, 69.73)
$row8 = @(9.37, 11.26, 13.23, 15.28, 17.21, 19.21, 21.28, 23.21, 25.15, 24.62, 22.05, 18.84, 18.76, 21.73, 24.52, 27.41, 30.41, 33.22, 36.05, 39.17, 42.01, 44.85, 47.81, 50.86, 53.72, 56.69, 59.76, 62.63, 65.14, 67.29, 69.23)
$row9 = @(17.92, 18.15, 18.87, 19.9, 21.02, 22.28, 23.66, 24.97, 26.27, 27.59, 28.41, 25.77, 17, 20.07, 23.59, 26.86, 30.06, 32.96, 35.8, 38.91, 41.69, 44.45, 47.29, 50.18, 52.86, 55.58, 58.34, 60.84, 63.24, 65.78, 67.94)
$row10 = @(12.95, 14.26, 15.62, 17.02, 18.3, 19.59, 20.8, 21.72, 22.25, 21.78, 19.11, 9.08, 12.74, 16.63, 21.02, 24.9, 28.57, 31.82, 34.96, 38.35, 41.37, 44.35, 47.41, 50.55, 53.47, 56.49, 59.61, 62.51, 65.41, 68.62, 71.52)

for ($i = 0; $i -lt $row7.Length; $i++) {
    $ws.Cells.Item(7, 2 + $i).Value = $row7[$i]
    $ws.Cells.Item(8, 2 + $i).Value = $row8[$i]
    $ws.Cells.Item(9, 2 + $i).Value = $row9[$i]
    $ws.Cells.Item(10, 2 + $i).Value = $row10[$i]
}

# Fix chart series formulas to reference the renamed sheet
$co = $ws.ChartObjects().Item(1)
$chart = $co.Chart
$s1 = $chart.SeriesCollection().Item(1)
$s1.Formula = "=SERIES(""Cremer"",'Madera_5.0_4.0_54.0'!`$B`$6:`$AF`$6,'Madera_5.0_4.0_54.0'!`$B`$7:`$AF`$7,1)"
$s2 = $chart.SeriesCollection().Item(2)
$s2.Formula = "=SERIES(""Sharp"",'Madera_5.0_4.0_54.0'!`$B`$6:`$AF`$6,'Madera_5.0_4.0_54.0'!`$B`$8:`$AF`$8,2)"
$s3 = $chart.SeriesCollection().Item(3)
$s3.Formula = "=SERIES(""Davy"",'Madera_5.0_4.0_54.0'!`$B`$6:`$AF`$6,'Madera_5.0_4.0_54.0'!`$B`$9:`$AF`$9,3)"
$s4 = $chart.SeriesCollection().Item(4)
$s4.Formula = "=SERIES('Madera_5.0_4.0_54.0'!`$A`$10,'Madera_5.0_4.0_54.0'!`$B`$6:`$AF`$6,'Madera_5.0_4.0_54.0'!`$B`$10:`$AF`$10,4)"

$wb.Save()
